$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.064.99'
$ws.Range('E2').Value = '  -2.17%  '
$ws.Range('D3').Value = '2.426.10'
$ws.Range('E3').Value = '  -1.08%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '572.46'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.89%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.33'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.29%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('E8').Value = '  -0.89%  '
$ws.Range('D9').Value = '2.412.22'
$ws.Range('E9').Value = '  -1.55%  '
$ws.Range('E10').Value = '  -0.79%  '
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.12'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.98%  '
$ws.Range('E13').Value = '  -1.92%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.22'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.10%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000172'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.37%  '
$ws.Range('D16').Value = '2.826.72'
$ws.Range('D17').Value = '60.921.58'
$ws.Range('E17').Value = '  -2.19%  '
$ws.Range('D18').Value = '2.410.96'
$ws.Range('E18').Value = '  -1.47%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.67'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +7.28%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.70'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '323.88'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.17%  '
$ws.Range('E22').Value = '  -1.54%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.09'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.47%  '
$ws.Range('E24').Value = '  +0.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.89'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.90%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '65.02'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.24%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '586.92'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.47%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.51'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -7.67%  '
$ws.Range('D29').Value = '2.541.20'
$ws.Range('E29').Value = '  -1.30%  '
$ws.Range('D30').Value = '0.0₃0942'
$ws.Range('E30').Value = '  -3.50%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.93'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.19%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.37'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.56%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.86'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.40%  '
$ws.Range('E34').Value = '  -0.93%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.66'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.70%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.42'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.39%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '152.04'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.33%  '
$ws.Range('E39').Value = '  -2.52%  '
$ws.Range('E40').Value = '  -0.60%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.17'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.71%  '
$ws.Range('E43').Value = '  -2.06%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.21'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.02%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.38'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.40%  '
$ws.Range('D46').Value = '0.0₆0289'
$ws.Range('E46').Value = '  +13.41%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '142.56'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.33%  '
$ws.Range('E48').Value = '  -3.55%  '
$ws.Range('E49').Value = '  -2.55%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.77'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.85%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0507'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.53%  '
